$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row: "<Feldname>_old" -> "<Feldname>_FV2310"
#    and "<Feldname>_new" -> "<Feldname>_FV2404" (the "diff" header is kept).
# ---------------------------------------------------------------------------
$newHeaders = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeaders[$i]
}

# ---------------------------------------------------------------------------
# 2) Turn the header + data range into an actual Excel Table ("Table1"),
#    ref A1:U80, with an AutoFilter on the header row.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U80"), 0, 1)
$lo.Name = "Table1"

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split after row 1, top-left of the scrollable
#    pane is A2).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
